# Update the cryptos list: refresh Price (column D) and Volume(1h) (column E)
# values for the rows scraped by the GitHub Actions job.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-Cell($row, $col, $value) {
    $cell = $ws.Cells.Item($row, $col)
    # Force text formatting so numeric-looking strings (e.g. "1.00", "173.00")
    # are preserved exactly as text rather than being auto-converted to numbers.
    $cell.NumberFormat = "@"
    $cell.Value = $value
}

# Row 2 - Bitcoin
Set-Cell 2 4 "66.463.79"
Set-Cell 2 5 "  -4.89%  "

# Row 3 - Ethereum
Set-Cell 3 4 "3.441.00"
Set-Cell 3 5 "  -6.75%  "

# Row 4 - TetherUSD
Set-Cell 4 4 "1.00"
Set-Cell 4 5 "  -0.01%  "

# Row 5 - BNB
Set-Cell 5 4 "599.60"
Set-Cell 5 5 "  -7.46%  "

# Row 6 - Solana
Set-Cell 6 4 "146.46"
Set-Cell 6 5 "  -9.62%  "

# Row 7 - LidoStakedEther
Set-Cell 7 4 "3.438.68"
Set-Cell 7 5 "  -6.77%  "

# Row 8 - USDC
Set-Cell 8 5 "  +0.06%  "

# Row 10 - Dogecoin
Set-Cell 10 5 "  -7.55%  "

# Row 11 - Toncoin
Set-Cell 11 5 "  -4.96%  "

# Row 12 - Cardano
Set-Cell 12 5 "  -6.49%  "

# Row 13 - ShibaInu
Set-Cell 13 4 "0.0000212"
Set-Cell 13 5 "  -8.78%  "

# Row 14 - WrappedliquidstakedEther2.0
Set-Cell 14 4 "4.027.74"
Set-Cell 14 5 "  -6.63%  "

# Row 15 - Avalanche
Set-Cell 15 4 "30.71"
Set-Cell 15 5 "  -6.25%  "

# Row 16 - WrappedEther
Set-Cell 16 4 "3.435.51"
Set-Cell 16 5 "  -6.66%  "

# Row 17 - WrappedBTC
Set-Cell 17 4 "66.401.22"
Set-Cell 17 5 "  -4.97%  "

# Row 18 - TRON
Set-Cell 18 5 "  -1.09%  "

# Row 19 - Polkadot
Set-Cell 19 5 "  -3.22%  "

# Row 20 - Chainlink
Set-Cell 20 4 "14.72"
Set-Cell 20 5 "  -7.87%  "

# Row 21 - BitcoinCash
Set-Cell 21 4 "436.34"
Set-Cell 21 5 "  -7.43%  "

# Row 22 - Uniswap
Set-Cell 22 4 "8.88"
Set-Cell 22 5 "  -14.13%  "

# Row 23 - Polygon
Set-Cell 23 4 "0.615"
Set-Cell 23 5 "  -5.72%  "

# Row 24 - Litecoin
Set-Cell 24 4 "76.47"
Set-Cell 24 5 "  -4.47%  "

# Row 25 - Dai
Set-Cell 25 5 "  +0.05%  "

# Row 26 - WrappedeETH
Set-Cell 26 4 "3.583.80"
Set-Cell 26 5 "  -6.60%  "

# Row 27 - PEPE
Set-Cell 27 5 "  -4.45%  "

# Row 28 - InternetComputer(DFINITY)
Set-Cell 28 5 "  -9.76%  "

# Row 29 - RenderToken
Set-Cell 29 4 "8.15"
Set-Cell 29 5 "  -10.81%  "

# Row 30 - PancakeSwap
Set-Cell 30 4 "2.49"
Set-Cell 30 5 "  -6.36%  "

# Row 31 - Binance-PegBSC-USD
Set-Cell 31 4 "0.999"
Set-Cell 31 5 "  -0.07%  "

# Row 32 - Fetch.AI
Set-Cell 32 4 "1.52"
Set-Cell 32 5 "  -11.50%  "

# Row 33 - Kaspa
Set-Cell 33 4 "0.158"
Set-Cell 33 5 "  -6.19%  "

# Row 34 - EthereumClassic
Set-Cell 34 4 "25.29"
Set-Cell 34 5 "  -5.46%  "

# Row 37 - RenzoRestakedETH
Set-Cell 37 4 "3.430.42"
Set-Cell 37 5 "  -6.96%  "

# Row 38 - Aptos
Set-Cell 38 5 "  -7.72%  "

# Row 39 - USDe
Set-Cell 39 5 "  +0.08%  "

# Row 41 - Monero
Set-Cell 41 4 "173.00"
Set-Cell 41 5 "  -4.22%  "

# Row 42 - Stacks
Set-Cell 42 5 "  -5.72%  "

# Row 43 - Hedera
Set-Cell 43 4 "0.0851"
Set-Cell 43 5 "  -5.74%  "

# Row 44 - Filecoin
Set-Cell 44 4 "5.34"
Set-Cell 44 5 "  -9.31%  "

# Row 45 - Mantle
Set-Cell 45 4 "0.871"
Set-Cell 45 5 "  -6.69%  "

# Row 46 - OKB
Set-Cell 46 4 "45.21"
Set-Cell 46 5 "  -2.94%  "

# Row 47 - ONDO
Set-Cell 47 4 "1.20"
Set-Cell 47 5 "  -4.39%  "

# Row 48 - InjectiveProtocol
Set-Cell 48 4 "25.76"
Set-Cell 48 5 "  -12.29%  "

# Row 49 - Cosmos
Set-Cell 49 5 "  -5.03%  "

# Row 50 - dogwifhat
Set-Cell 50 4 "2.45"
Set-Cell 50 5 "  -15.06%  "

# Row 51 - SuiNetwork
Set-Cell 51 4 "0.987"
Set-Cell 51 5 "  -6.75%  "
